$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '54.386.62'
$ws.Range('E2').Value = '  -6.79%  '

# Row 3
$ws.Range('D3').Value = '2.877.50'
$ws.Range('E3').Value = '  -9.83%  '

# Row 4
$ws.Range('E4').Value = '  -0.12%  '

# Row 5
$ws.Range('D5').Value = '''470.84'
$ws.Range('E5').Value = '  -11.64%  '

# Row 6
$ws.Range('D6').Value = '''125.96'
$ws.Range('E6').Value = '  -6.58%  '

# Row 7
$ws.Range('E7').Value = '  -0.07%  '

# Row 8
$ws.Range('D8').Value = '2.870.38'
$ws.Range('E8').Value = '  -10.10%  '

# Row 9
$ws.Range('D9').Value = '''0.403'
$ws.Range('E9').Value = '  -11.36%  '

# Row 10
$ws.Range('D10').Value = '''6.62'
$ws.Range('E10').Value = '  -9.71%  '

# Row 11
$ws.Range('D11').Value = '''0.0959'
$ws.Range('E11').Value = '  -14.56%  '

# Row 12
$ws.Range('D12').Value = '''0.328'
$ws.Range('E12').Value = '  -16.79%  '

# Row 13
$ws.Range('E13').Value = '  -4.91%  '

# Row 14
$ws.Range('D14').Value = '3.360.33'
$ws.Range('E14').Value = '  -10.05%  '

# Row 15
$ws.Range('D15').Value = '''23.16'
$ws.Range('E15').Value = '  -9.84%  '

# Row 16
$ws.Range('D16').Value = '54.342.81'
$ws.Range('E16').Value = '  -7.08%  '

# Row 17
$ws.Range('D17').Value = '2.876.60'
$ws.Range('E17').Value = '  -9.66%  '

# Row 18
$ws.Range('D18').Value = '''0.0000133'
$ws.Range('E18').Value = '  -14.48%  '

# Row 19
$ws.Range('D19').Value = '''5.32'
$ws.Range('E19').Value = '  -9.31%  '

# Row 20
$ws.Range('D20').Value = '''11.39'
$ws.Range('E20').Value = '  -13.96%  '

# Row 21
$ws.Range('D21').Value = '''7.05'
$ws.Range('E21').Value = '  -12.93%  '

# Row 22
$ws.Range('D22').Value = '''297.62'
$ws.Range('E22').Value = '  -17.13%  '

# Row 23
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.26%  '

# Row 24
$ws.Range('D24').Value = '''0.439'
$ws.Range('E24').Value = '  -14.89%  '

# Row 25
$ws.Range('D25').Value = '''58.51'
$ws.Range('E25').Value = '  -16.01%  '

# Row 26
$ws.Range('E26').Value = '  +0.01%  '

# Row 27
$ws.Range('E27').Value = '  -10.13%  '

# Row 28
$ws.Range('E28').Value = '  -0.13%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0806'
$ws.Range('E29').Value = '  -15.02%  '

# Row 30
$ws.Range('D30').Value = '''6.10'
$ws.Range('E30').Value = '  -12.38%  '

# Row 31
$ws.Range('D31').Value = '''1.12'
$ws.Range('E31').Value = '  -7.19%  '

# Row 32
$ws.Range('D32').Value = '''6.14'
$ws.Range('E32').Value = '  -12.14%  '

# Row 33
$ws.Range('D33').Value = '''1.61'
$ws.Range('E33').Value = '  -15.56%  '

# Row 34
$ws.Range('D34').Value = '''18.56'
$ws.Range('E34').Value = '  -14.36%  '

# Row 35
$ws.Range('D35').Value = '''4.16'
$ws.Range('E35').Value = '  -15.81%  '

# Row 36
$ws.Range('D36').Value = '''137.22'
$ws.Range('E36').Value = '  -14.62%  '

# Row 37
$ws.Range('D37').Value = '''5.38'
$ws.Range('E37').Value = '  -14.67%  '

# Row 38
$ws.Range('E38').Value = '  -15.00%  '

# Row 39
$ws.Range('D39').Value = '''22.95'
$ws.Range('E39').Value = '  -10.86%  '

# Row 40
$ws.Range('B40').Value = 'RenzoRestakedETH'
$ws.Range('C40').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D40').Value = '2.894.82'
$ws.Range('E40').Value = '  -10.06%  '

# Row 41
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '''0.0614'
$ws.Range('E41').Value = '  -12.88%  '

# Row 42
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.05%  '

# Row 43
$ws.Range('D43').Value = '''34.89'
$ws.Range('E43').Value = '  -14.35%  '

# Row 44
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '''0.601'
$ws.Range('E44').Value = '  -15.04%  '

# Row 45
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').Value = '''0.942'
$ws.Range('E45').Value = '  -13.15%  '

# Row 46
$ws.Range('D46').Value = '''1.31'
$ws.Range('E46').Value = '  -11.36%  '

# Row 47
$ws.Range('D47').Value = '''3.37'
$ws.Range('E47').Value = '  -15.88%  '

# Row 48
$ws.Range('D48').Value = '2.039.21'
$ws.Range('E48').Value = '  -10.88%  '

# Row 49
$ws.Range('D49').Value = '''5.33'
$ws.Range('E49').Value = '  -14.41%  '

# Row 50
$ws.Range('D50').Value = '''17.77'
$ws.Range('E50').Value = '  -13.28%  '

# Row 51
$ws.Range('D51').Value = '''0.0213'
$ws.Range('E51').Value = '  -10.55%  '
